# Update "想去人数" (F column) counts on both the "展览" and "全部类型"
# worksheets, which carry identical data in this workbook.

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 1076
    3  = 777
    5  = 35
    6  = 1103
    8  = 1911
    9  = 6781
    11 = 387
    12 = 321
    13 = 113
    14 = 384
    15 = 144
    16 = 6998
    17 = 283
    18 = 1311
    19 = 141
    21 = 224
    22 = 120
    23 = 286
    24 = 120
    27 = 104
    28 = 14
    29 = 399
    30 = 599
    32 = 84
    34 = 65
    35 = 29
    36 = 67
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
